$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "_old" / "_new" header-name suffixes to the respective
#    format-version suffixes ("_FV2310" / "_FV2404").
$headers = @{
    "A1" = "Segmentname_FV2310";
    "B1" = "Segmentgruppe_FV2310";
    "C1" = "Segment_FV2310";
    "D1" = "Datenelement_FV2310";
    "E1" = "Segment ID_FV2310";
    "F1" = "Code_FV2310";
    "G1" = "Qualifier_FV2310";
    "H1" = "Beschreibung_FV2310";
    "I1" = "Bedingungsausdruck_FV2310";
    "J1" = "Bedingung_FV2310";
    "K1" = "diff";
    "L1" = "Segmentname_FV2404";
    "M1" = "Segmentgruppe_FV2404";
    "N1" = "Segment_FV2404";
    "O1" = "Datenelement_FV2404";
    "P1" = "Segment ID_FV2404";
    "Q1" = "Code_FV2404";
    "R1" = "Qualifier_FV2404";
    "S1" = "Beschreibung_FV2404";
    "T1" = "Bedingungsausdruck_FV2404";
    "U1" = "Bedingung_FV2404";
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# 2) Turn the header row + data range into a real Excel Table ("Table1"),
#    matching the header names that are now in row 1.
$rng = $ws.Range("A1:U60")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# 3) Freeze the header row (split/freeze after row 1) and select the pane
#    below the freeze, mirroring the sheetView pane/selection added in the
#    diff.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
